$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

$find.Execute(
    "Kampagnendaten 2022 für das Sternbild Perseus-Konstellation",
    $true,
    $false,
    $false,
    $false,
    $false,
    $true,
    1,
    $false,
    "Kampagnendaten 2022 für das Perseus-Konstellation",
    2
)
